# Penalty Reward System (unfinished) - forecast_summary_B09J64TBJG.xlsx
#
# This shifts the "Forecast Comparison" week-start dates forward by one
# week and zeroes out most of the MyForecast (column D) values (an
# in-progress penalty/reward experiment), and pokes a handful of
# (now-stale/inconsistent) numbers on the "Summary" sheet to match.
#
# Helper: write a value as LITERAL TEXT (mirrors the original inlineStr
# cells) without leaving the cell tagged with a different number format /
# style than it started with - Excel's normal .Value setter auto-detects
# dates / numbers from date- or number-shaped strings, which would change
# both the stored type and the cell style.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Forecast Comparison
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# row -> (new Week_Start_Date, new MyForecast)
Set-TextValue $ws1.Range("B2") "2025-01-12"
$ws1.Range("D2").Value = 0

Set-TextValue $ws1.Range("B3") "2025-01-19"
$ws1.Range("D3").Value = 0

Set-TextValue $ws1.Range("B4") "2025-01-26"
$ws1.Range("D4").Value = 0

Set-TextValue $ws1.Range("B5") "2025-02-02"
$ws1.Range("D5").Value = 0

Set-TextValue $ws1.Range("B6") "2025-02-09"
$ws1.Range("D6").Value = 0

Set-TextValue $ws1.Range("B7") "2025-02-16"
$ws1.Range("D7").Value = 0

Set-TextValue $ws1.Range("B8") "2025-02-23"
$ws1.Range("D8").Value = 1

Set-TextValue $ws1.Range("B9") "2025-03-02"
$ws1.Range("D9").Value = 0

Set-TextValue $ws1.Range("B10") "2025-03-09"
$ws1.Range("D10").Value = 0

Set-TextValue $ws1.Range("B11") "2025-03-16"
$ws1.Range("D11").Value = 0

Set-TextValue $ws1.Range("B12") "2025-03-23"
$ws1.Range("D12").Value = 0

Set-TextValue $ws1.Range("B13") "2025-03-30"
$ws1.Range("D13").Value = 0

Set-TextValue $ws1.Range("B14") "2025-04-06"
$ws1.Range("D14").Value = 0

Set-TextValue $ws1.Range("B15") "2025-04-13"
$ws1.Range("D15").Value = 0

Set-TextValue $ws1.Range("B16") "2025-04-20"
$ws1.Range("D16").Value = 0

Set-TextValue $ws1.Range("B17") "2025-04-27"
$ws1.Range("D17").Value = 0

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

Set-TextValue $ws2.Range("B2") "2023-01-15 to 2025-01-05"
Set-TextValue $ws2.Range("B9") "6"
Set-TextValue $ws2.Range("B10") "3"
Set-TextValue $ws2.Range("B11") "2"
Set-TextValue $ws2.Range("B12") "1"
Set-TextValue $ws2.Range("B13") "2025-02-23"
Set-TextValue $ws2.Range("B14") "0"
Set-TextValue $ws2.Range("B15") "2025-03-02"
